$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 32 (pushes existing rows 32-45 down to 33-46,
# carrying their values/styles/heights with them).
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the Gini script entry.
$ws.Range("A32").Value = "TCR_analyses/giniscript.py"
$ws.Range("B32").Value = "Calculate Gini coefficients for TCR clonotype data"
$ws.Range("C32").Value = "Text files containing clonotype frequency data, from TCR_circos.R script"
$ws.Range("D32").Value = "Text files containing Gini TCR coefficients"

# The row insert does not automatically relocate the worksheet's hyperlinks,
# so the two existing hyperlinks (previously anchored at A35/A36) still point
# at those old cell addresses even though their text moved down to A36/A37.
# Remove the stale hyperlinks and re-add them at the correct, shifted cells.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A36"), "https://github.com/IzarLab/Melanoma_Brain_Metastasis/Spatial_analysis/rctd_mbpm_puckdata.R", "", "", "https://github.com/IzarLab/Melanoma_Brain_Metastasis/Spatial_analysis/rctd_mbpm_puckdata.R")
$ws.Hyperlinks.Add($ws.Range("A37"), "https://github.com/IzarLab/Melanoma_Brain_Metastasis/Spatial_analysis/rctd_mbpm_puckdata_plot_cell_types_pub_quality.R", "", "", "https://github.com/IzarLab/Melanoma_Brain_Metastasis/Spatial_analysis/rctd_mbpm_puckdata_plot_cell_types_pub_quality.R")

# Match the other incidental worksheet metadata changes from the edit:
# the selected cell and the sheet's default column width.
$ws.Range("B31").Select()
$ws.StandardWidth = 11.625
